$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '297.99'
$ws.Range('E2').Value = '1.91%'
$ws.Range('G2').Value = '7'
$ws.Range('D3').Value = '41.77'
$ws.Range('E3').Value = '3.12%'
$ws.Range('G3').Value = '7'
$ws.Range('D4').Value = '5.021'
$ws.Range('E4').Value = '-0.33%'
$ws.Range('G4').Value = '7'
$ws.Range('D5').Value = '0.07535'
$ws.Range('E5').Value = '2.83%'
$ws.Range('G5').Value = '7'
$ws.Range('B6').Value = 'FTXToken'
$ws.Range('C6').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D6').Value = '1.624'
$ws.Range('E6').Value = '6.24%'
$ws.Range('G6').Value = '7'
$ws.Range('B7').Value = 'MXToken'
$ws.Range('C7').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D7').Value = '0.9216'
$ws.Range('E7').Value = '-0.96%'
$ws.Range('G7').Value = '7'
$ws.Range('B8').Value = 'BTSEToken'
$ws.Range('C8').Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range('D8').Value = '2.401'
$ws.Range('E8').Value = '1.83%'
$ws.Range('G8').Value = '7'
$ws.Range('B9').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C9').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('D9').Value = '0.1182'
$ws.Range('E9').Value = '0.65%'
$ws.Range('G9').Value = '7'
$ws.Range('B10').Value = 'WazirX'
$ws.Range('C10').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('D10').Value = '0.1827'
$ws.Range('E10').Value = '5.12%'
$ws.Range('G10').Value = '7'
$ws.Range('B11').Value = 'MandalaExchangeToken'
$ws.Range('C11').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D11').Value = '0.08952'
$ws.Range('E11').Value = '3.11%'
$ws.Range('G11').Value = '7'
$ws.Range('B12').Value = 'BitrueCoin'
$ws.Range('C12').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D12').Value = '0.04086'
$ws.Range('E12').Value = '-5.82%'
$ws.Range('G12').Value = '7'
$ws.Range('B13').Value = 'BitMartToken'
$ws.Range('C13').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D13').Value = '0.1050'
$ws.Range('E13').Value = '-0.40%'
$ws.Range('G13').Value = '7'
$ws.Range('B14').Value = 'BitForexToken'
$ws.Range('C14').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D14').Value = '0.001280'
$ws.Range('E14').Value = '0.53%'
$ws.Range('G14').Value = '7'
$ws.Range('B15').Value = 'TigerCash'
$ws.Range('C15').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('D15').Value = '0.005788'
$ws.Range('E15').Value = '-3.49%'
$ws.Range('G15').Value = '7'
$ws.Range('B16').Value = 'LEO'
$ws.Range('C16').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D16').Value = '3.340'
$ws.Range('E16').Value = '0.12%'
$ws.Range('G16').Value = '7'
$ws.Range('B17').Value = 'GateToken'
$ws.Range('C17').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range('D17').Value = '4.377'
$ws.Range('E17').Value = '2.26%'
$ws.Range('G17').Value = '7'
$ws.Range('D18').Value = '0.3328'
$ws.Range('E18').Value = '1.14%'
$ws.Range('G18').Value = '7'
$ws.Range('D19').Value = '8.274'
$ws.Range('E19').Value = '3.79%'
$ws.Range('G19').Value = '7'
$ws.Range('E20').Value = '-2.06%'
$ws.Range('G20').Value = '7'
$ws.Range('D21').Value = '0.3222'
$ws.Range('G21').Value = '7'
$ws.Range('D22').Value = '0.04075'
$ws.Range('E22').Value = '3.29%'
$ws.Range('G22').Value = '7'
$ws.Range('E23').Value = '0.34%'
$ws.Range('G23').Value = '7'
$ws.Range('D24').Value = '0.003893'
$ws.Range('E24').Value = '2.85%'
$ws.Range('G24').Value = '7'
$ws.Range('E25').Value = '-3.88%'
$ws.Range('G25').Value = '7'
$ws.Range('G26').Value = '7'
$ws.Range('G27').Value = '7'
$ws.Range('G28').Value = '7'
$ws.Range('G29').Value = '7'
$ws.Range('G30').Value = '7'
$ws.Range('G31').Value = '7'
$ws.Range('G32').Value = '7'
$ws.Range('G33').Value = '7'
$ws.Range('G34').Value = '7'
$ws.Range('G35').Value = '7'
$ws.Range('G36').Value = '7'
$ws.Range('G37').Value = '7'
$ws.Range('D38').Value = '0.02406'
$ws.Range('E38').Value = '5.10%'
$ws.Range('G38').Value = '7'
$ws.Range('D39').Value = '0.05215'
$ws.Range('E39').Value = '3.45%'
$ws.Range('G39').Value = '7'
$ws.Range('E40').Value = '1.21%'
$ws.Range('G40').Value = '7'
$ws.Range('D41').Value = '0.007824'
$ws.Range('E41').Value = '1.12%'
$ws.Range('G41').Value = '7'
$ws.Range('E42').Value = '2.97%'
$ws.Range('G42').Value = '7'
$ws.Range('D43').Value = '0.007393'
$ws.Range('E43').Value = '0.44%'
$ws.Range('G43').Value = '7'
$ws.Range('D44').Value = '0.007619'
$ws.Range('E44').Value = '-7.88%'
$ws.Range('G44').Value = '7'
$ws.Range('D45').Value = '0.2968'
$ws.Range('E45').Value = '1.63%'
$ws.Range('G45').Value = '7'
$ws.Range('D46').Value = '0.00006575'
$ws.Range('E46').Value = '4.61%'
$ws.Range('G46').Value = '7'
$ws.Range('E47').Value = '-0.03%'
$ws.Range('G47').Value = '7'
$ws.Range('D48').Value = '0.04746'
$ws.Range('E48').Value = '45.74%'
$ws.Range('G48').Value = '7'
$ws.Range('G49').Value = '7'
$ws.Range('E50').Value = '-0.03%'
$ws.Range('G50').Value = '7'
$ws.Range('D51').Value = '0.0002001'
$ws.Range('E51').Value = '-0.03%'
$ws.Range('G51').Value = '7'
